$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.974.78'
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').Value = '2.385.64'
$ws.Range('E3').Value = '  +3.53%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '300.85'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '99.08'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.76%  '
$ws.Range('E7').Value = '  -0.91%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.510'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.52%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.54'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -6.23%  '
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  -3.39%  '
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('D14').Value = '2.749.38'
$ws.Range('E14').Value = '  +3.43%  '
$ws.Range('D15').Value = '2.363.82'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.820'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.73'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('D18').Value = '45.906.15'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('E19').Value = '  -3.60%  '
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.08'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.05'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '244.21'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.57%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.81'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.71%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.93'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.69%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '39.63'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -9.20%  '
$ws.Range('E28').Value = '  -3.21%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.77'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.40%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.82'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +20.42%  '
$ws.Range('E31').Value = '  +4.84%  '
$ws.Range('E32').Value = '  +7.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '146.90'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0773'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.62%  '
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('E37').Value = '  +6.82%  '
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '14.89'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.57%  '
$ws.Range('E40').Value = '  -3.53%  '
$ws.Range('E41').Value = '  -1.80%  '
$ws.Range('E42').Value = '  -6.82%  '
$ws.Range('D43').Value = '1.939.81'
$ws.Range('E43').Value = '  +4.63%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '91.89'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +4.08%  '
$ws.Range('E46').Value = '  -9.73%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.48'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +5.37%  '
$ws.Range('E48').Value = '  -5.17%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '99.01'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('D50').Value = '2.619.83'
$ws.Range('E50').Value = '  +3.42%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '68.59'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -8.04%  '
